$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting dbExcel/WebData columns right
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Range("B1").Value = "StatQuery"

# New query text for the inserted column (row 2), matching the wrap style used in A2
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Serous endometrial adenocarcinoma']   OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").WrapText = $true

# Column width for the newly inserted column (match column A's width as closely as possible)
$ws.Columns.Item(2).ColumnWidth = 75

# Move selection
$ws.Range("A2").Select()
